$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the disclaimer text date from 2021-03-26 to 2021-03-29
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-29 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values
$ws.Range("D2").Value = 0.2550365200415783
$ws.Range("E2").Value = -0.003388128999874507

$ws.Range("D3").Value = 0.25079923501389
$ws.Range("E3").Value = -0.008728542333430234

$ws.Range("D4").Value = 0.2517523377561418
$ws.Range("E4").Value = -0.005200874349890805

$ws.Range("D5").Value = 0.24241190718839
$ws.Range("E5").Value = 0.01023654724028211

$ws.Range("E6").Value = -0.001881079706010591

# Restore sheet protection (sheet was protected before the edit)
$ws.Protect()
